$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 296.7258181216166
$ws.Range("P2").Value = 785.5721091280598
$ws.Range("Q2").Value = 43.96014014652901
$ws.Range("R2").Value = 741.6119689815307
$ws.Range("S2").Value = 1200.359498886298
$ws.Range("T2").Value = -458.7475299047676
$ws.Range("U2").Value = 596.8566180376354
$ws.Range("V2").Value = 33.39973539213398
$ws.Range("W2").Value = 563.4568826455014
$ws.Range("X2").Value = 912.0009514210426
$ws.Range("Y2").Value = -348.5440687755411
$ws.Range("Z2").Value = 710.517092521723
$ws.Range("AA2").Value = 26.00149233543534
$ws.Range("AB2").Value = 684.5156001862878
$ws.Range("AC2").Value = 862.133577495752
$ws.Range("AD2").Value = -177.6179773094641
$ws.Range("AE2").Value = -208.9144135003577
$ws.Range("AF2").Value = 0.6673262650571995
$ws.Range("AG2").Value = -209.5817397654149
$ws.Range("AH2").Value = -336.2788336308169
$ws.Range("AI2").Value = 126.697093865402
$ws.Range("AJ2").Value = -3066.442921159479
$ws.Range("AK2").Value = -304.3761421095356
$ws.Range("AL2").Value = -2762.066779049943
$ws.Range("AM2").Value = -2532.8262049561
$ws.Range("AN2").Value = -229.2405740938434
$ws.Range("AO2").Value = 576.657695627702
$ws.Range("AP2").Value = 44.62746641158621
$ws.Range("AQ2").Value = 532.0302292161158
$ws.Range("AR2").Value = 859.9252386767228
$ws.Range("AS2").Value = -327.895009460607
$ws.Range("AT2").Value = 387.9422045372778
$ws.Range("AU2").Value = 34.06706165719118
$ws.Range("AV2").Value = 353.8751428800866
$ws.Range("AW2").Value = 572.564936131811
$ws.Range("AX2").Value = -218.6897932517244
$ws.Range("AY2").Value = 501.6026790213653
$ws.Range("AZ2").Value = 26.66881860049255
$ws.Range("BA2").Value = 474.9338604208727
$ws.Range("BB2").Value = 524.3654869153245
$ws.Range("BC2").Value = -49.43162649445171
$ws.Range("BD2").Value = -3275.357334659836
$ws.Range("BE2").Value = -303.7088158444784
$ws.Range("BF2").Value = -2971.648518815358
$ws.Range("BG2").Value = -2841.589093353845
$ws.Range("BH2").Value = -130.0594254615129
$ws.Range("BI2").Value = -2489.785225531777
$ws.Range("BJ2").Value = -259.7486756979494
$ws.Range("BK2").Value = -2230.036549833827
$ws.Range("BL2").Value = -1645.384997982764
$ws.Range("BM2").Value = -584.6515518510636
$ws.Range("BN2").Value = -2678.500716622201
$ws.Range("BO2").Value = -270.3090804523443
$ws.Range("BP2").Value = -2408.191636169857
$ws.Range("BQ2").Value = -1932.745306068156
$ws.Range("BR2").Value = -475.4463301017009
$ws.Range("BS2").Value = -2564.840242138113
$ws.Range("BT2").Value = -277.707323509043
$ws.Range("BU2").Value = -2287.13291862907
$ws.Range("BV2").Value = -1980.945312729317
$ws.Range("BW2").Value = -306.1876058997539
$ws.Range("C3").Value = 67.87593501224863
$ws.Range("P3").Value = 412.131421676704
$ws.Range("Q3").Value = 2.717596839911276
$ws.Range("R3").Value = 399.6376714038686
$ws.Range("S3").Value = 281.9201185032147
$ws.Range("T3").Value = 117.7175529006539
$ws.Range("U3").Value = 313.1264000729696
$ws.Range("V3").Value = 2.064757188056933
$ws.Range("W3").Value = 303.6339837208563
$ws.Range("X3").Value = 214.1953444265777
$ws.Range("Y3").Value = 89.4386392942786
$ws.Range("Z3").Value = 259.2630175815141
$ws.Range("AA3").Value = 1.78631582530416
$ws.Range("AB3").Value = 252.6284294332987
$ws.Range("AC3").Value = 201.3742944093364
$ws.Range("AD3").Value = 51.25413502396236
$ws.Range("AE3").Value = -214.8798458902175
$ws.Range("AF3").Value = -1.753310944870222
$ws.Range("AG3").Value = -211.5516686820544
$ws.Range("AH3").Value = -75.01771291657376
$ws.Range("AI3").Value = -136.5339557654807
$ws.Range("AJ3").Value = -34.80573054928919
$ws.Range("AK3").Value = -0.1252237217194288
$ws.Range("AL3").Value = -34.2763366601559
$ws.Range("AM3").Value = -648.8838963646766
$ws.Range("AN3").Value = 614.6075597045208
$ws.Range("AO3").Value = 197.251575786486
$ws.Range("AP3").Value = 0.9642858950410502
$ws.Range("AQ3").Value = 188.0860027218137
$ws.Range("AR3").Value = 205.9518530995217
$ws.Range("AS3").Value = -17.86585037770804
$ws.Range("AT3").Value = 98.24655418275236
$ws.Range("AU3").Value = 0.311446243186712
$ws.Range("AV3").Value = 92.08231503880211
$ws.Range("AW3").Value = 138.4554272258436
$ws.Range("AX3").Value = -46.37311218704151
$ws.Range("AY3").Value = 44.38317169129616
$ws.Range("AZ3").Value = 0.03300488043393627
$ws.Range("BA3").Value = 41.0767607512438
$ws.Range("BB3").Value = 126.0159144467731
$ws.Range("BC3").Value = -84.9391536955293
$ws.Range("BD3").Value = -249.6855764395069
$ws.Range("BE3").Value = -1.878534666589652
$ws.Range("BF3").Value = -245.8280053422106
$ws.Range("BG3").Value = -717.6073457006034
$ws.Range("BH3").Value = 471.7793403583929
$ws.Range("BI3").Value = 162.4458452371955
$ws.Range("BJ3").Value = 0.8390621733216187
$ws.Range("BK3").Value = 153.8096660616565
$ws.Range("BL3").Value = -436.6377744087304
$ws.Range("BM3").Value = 590.4474404703869
$ws.Range("BN3").Value = 63.44082363346268
$ws.Range("BO3").Value = 0.1862225214672813
$ws.Range("BP3").Value = 57.80597837864586
$ws.Range("BQ3").Value = -504.1342015497914
$ws.Range("BR3").Value = 561.9401799284373
$ws.Range("BS3").Value = 9.577441142006318
$ws.Range("BT3").Value = -0.09221884128549353
$ws.Range("BU3").Value = 6.800424091087336
$ws.Range("BV3").Value = -516.5738418441512
$ws.Range("BW3").Value = 523.3742659352384
$ws.Range("C4").Value = 40.27682264769349
$ws.Range("P4").Value = 672.7433163074238
$ws.Range("Q4").Value = 10.15962064363052
$ws.Range("R4").Value = 662.5836956637934
$ws.Range("S4").Value = 158.7412243708941
$ws.Range("T4").Value = 503.8424712928992
$ws.Range("U4").Value = 511.132327526684
$ws.Range("V4").Value = 7.719007265460548
$ws.Range("W4").Value = 503.4133202612235
$ws.Range("X4").Value = 120.6073245476186
$ws.Range("Y4").Value = 382.8059957136049
$ws.Range("Z4").Value = 501.2305127401583
$ws.Range("AA4").Value = 5.523434251378883
$ws.Range("AB4").Value = 495.7070784887794
$ws.Range("AC4").Value = 115.0298055079141
$ws.Range("AD4").Value = 380.6772729808653
$ws.Range("AE4").Value = -15.46973475103519
$ws.Range("AF4").Value = 0.1361768030811918
$ws.Range("AG4").Value = -15.60591155411638
$ws.Range("AH4").Value = -45.6912459449094
$ws.Range("AI4").Value = 30.08533439079302
$ws.Range("AJ4").Value = -2213.035655791258
$ws.Range("AK4").Value = -43.52422168775004
$ws.Range("AL4").Value = -2169.511434103508
$ws.Range("AM4").Value = -341.5905501638384
$ws.Range("AN4").Value = -1827.920883939669
$ws.Range("AO4").Value = 657.2735815563885
$ws.Range("AP4").Value = 10.29579744671172
$ws.Range("AQ4").Value = 646.9777841096768
$ws.Range("AR4").Value = 112.4859311839965
$ws.Range("AS4").Value = 534.4918529256803
$ws.Range("AT4").Value = 495.6625927756492
$ws.Range("AU4").Value = 7.855184068541742
$ws.Range("AV4").Value = 487.8074087071074
$ws.Range("AW4").Value = 74.48753063702324
$ws.Range("AX4").Value = 413.3198780700841
$ws.Range("AY4").Value = 485.7607779891227
$ws.Range("AZ4").Value = 5.659611054460074
$ws.Range("BA4").Value = 480.1011669346626
$ws.Range("BB4").Value = 69.13641153821867
$ws.Range("BC4").Value = 410.9647553964439
$ws.Range("BD4").Value = -2228.505390542293
$ws.Range("BE4").Value = -43.38804488466886
$ws.Range("BF4").Value = -2185.117345657624
$ws.Range("BG4").Value = -383.5468503155179
$ws.Range("BH4").Value = -1801.570495342107
$ws.Range("BI4").Value = -1555.762074234869
$ws.Range("BJ4").Value = -33.22842424103832
$ws.Range("BK4").Value = -1522.533649993831
$ws.Range("BL4").Value = -225.3696700560244
$ws.Range("BM4").Value = -1297.163979937807
$ws.Range("BN4").Value = -1717.373063015608
$ws.Range("BO4").Value = -35.6690376192083
$ws.Range("BP4").Value = -1681.7040253964
$ws.Range("BQ4").Value = -263.3680713550485
$ws.Range("BR4").Value = -1418.335954041351
$ws.Range("BS4").Value = -1727.274877802135
$ws.Range("BT4").Value = -37.86461063328997
$ws.Range("BU4").Value = -1689.410267168845
$ws.Range("BV4").Value = -268.7192661200032
$ws.Range("BW4").Value = -1420.691001048842
$ws.Range("C5").Value = 257.4628107574233
$ws.Range("P5").Value = 897.9706516325189
$ws.Range("Q5").Value = 16.74402402343943
$ws.Range("R5").Value = 863.5032577139413
$ws.Range("S5").Value = 1062.926204115415
$ws.Range("T5").Value = -199.4229464014741
$ws.Range("U5").Value = 682.2540158984511
$ws.Range("V5").Value = 12.72166034772225
$ws.Range("W5").Value = 656.066614477534
$ws.Range("X5").Value = 807.5828202659474
$ws.Range("Y5").Value = -151.5162057884133
$ws.Range("Z5").Value = 508.2468162466416
$ws.Range("AA5").Value = 5.566427996696711
$ws.Range("AB5").Value = 493.3409793565668
$ws.Range("AC5").Value = 765.0498697642827
$ws.Range("AD5").Value = -271.7088904077158
$ws.Range("AE5").Value = -309.9218189660585
$ws.Range("AF5").Value = -7.919478470079362
$ws.Range("AG5").Value = -286.6782318142385
$ws.Range("AH5").Value = -283.2836347470828
$ws.Range("AI5").Value = -3.394597067155728
$ws.Range("AJ5").Value = -957.7755509024936
$ws.Range("AK5").Value = -14.62481029835312
$ws.Range("AL5").Value = -915.9705474905128
$ws.Range("AM5").Value = -2447.156559409366
$ws.Range("AN5").Value = 1531.186011918853
$ws.Range("AO5").Value = 588.0488326664595
$ws.Range("AP5").Value = 8.824545553360059
$ws.Range("AQ5").Value = 576.8250258997019
$ws.Range("AR5").Value = 776.0369922950736
$ws.Range("AS5").Value = -199.2119663953717
$ws.Range("AT5").Value = 372.3321969323927
$ws.Range("AU5").Value = 4.802181877642883
$ws.Range("AV5").Value = 369.3883826632954
$ws.Range("AW5").Value = 521.5597647667612
$ws.Range("AX5").Value = -152.1713821034658
$ws.Range("AY5").Value = 198.324997280582
$ws.Range("AZ5").Value = -2.353050473382658
$ws.Range("BA5").Value = 206.6627475423271
$ws.Range("BB5").Value = 480.4740377874003
$ws.Range("BC5").Value = -273.8112902450732
$ws.Range("BD5").Value = -1267.697369868552
$ws.Range("BE5").Value = -22.54428876843248
$ws.Range("BF5").Value = -1202.648779304752
$ws.Range("BG5").Value = -2706.565181711626
$ws.Range("BH5").Value = 1503.916402406875
$ws.Range("BI5").Value = -369.726718236036
$ws.Range("BJ5").Value = -5.800264744993056
$ws.Range("BK5").Value = -339.1455215908126
$ws.Range("BL5").Value = -1647.244534657715
$ws.Range("BM5").Value = 1308.099013066903
$ws.Range("BN5").Value = -585.4433539701009
$ws.Range("BO5").Value = -9.822628420710227
$ws.Range("BP5").Value = -546.5821648272172
$ws.Range("BQ5").Value = -1901.721766993386
$ws.Range("BR5").Value = 1355.139602166169
$ws.Range("BS5").Value = -759.4505536219121
$ws.Range("BT5").Value = -16.97786077173577
$ws.Range("BU5").Value = -709.307799948186
$ws.Range("BV5").Value = -1942.807977655878
$ws.Range("BW5").Value = 1233.500177707692
